# edit.ps1 - apply "update logo, fix broken links" changes
#
# 1) Slide 1 ("logo") group shape: rebrand the two caption textboxes
#      "0030" -> "TIED"   (and recolor 1ABC9C -> 6D8E79)
#      "GEOG" -> "SA"
# 2) Refresh the cached datetimeFigureOut field text (10/13/23 -> 8/8/24)
#    on the Slide Master and every Custom Layout's Date placeholder.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Slide 1 group shape text + color updates
# ---------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)

for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $shp = $group.GroupItems.Item($i)

    if ($shp.Name -eq "TextBox 3") {
        # big "0030" number -> "TIED", recolor teal -> sage green
        $shp.TextFrame.TextRange.Text = "TIED"
        $shp.TextFrame.TextRange.Font.Color.RGB = 0x798E6D  # RGB(0x6D,0x8E,0x79) == 6D8E79
    }
    elseif ($shp.Name -eq "TextBox 7") {
        # "GEOG" -> "SA" (color unchanged)
        $shp.TextFrame.TextRange.Text = "SA"
    }
}

# ---------------------------------------------------------------
# 2) Re-cache the datetimeFigureOut field text everywhere it appears
# ---------------------------------------------------------------
$newDate = "8/8/24"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "10/13/23") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "10/13/23") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
